# Weekly update: insert a new price-observation row at the top of the data
# block (row 313), pushing the existing rows 313:389 down to 314:390.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 313 (shifts 313:389 -> 314:390)
$ws.Rows.Item(313).Insert()

# Populate the new row 313 with the new weekly observation.
$ws.Cells.Item(313, 1).Value = 10
$ws.Cells.Item(313, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(313, 3).Value = "La Araucanía"
$ws.Cells.Item(313, 4).Value = 44642
$ws.Cells.Item(313, 5).Value = 9
$ws.Cells.Item(313, 6).Value = "Fruta"
$ws.Cells.Item(313, 7).Value = 100101
$ws.Cells.Item(313, 8).Value = "Berries"
$ws.Cells.Item(313, 9).Value = 100101007
$ws.Cells.Item(313, 10).Value = "Kiwi"
$ws.Cells.Item(313, 11).Value = "Gold"
$ws.Cells.Item(313, 12).Value = "Segunda"
$ws.Cells.Item(313, 13).Value = 30
$ws.Cells.Item(313, 14).Value = 12000
$ws.Cells.Item(313, 15).Value = 12000
$ws.Cells.Item(313, 16).Value = 12000
$ws.Cells.Item(313, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(313, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(313, 19).Value = 667
$ws.Cells.Item(313, 20).Value = 18
